# ModelRuns_RTP2025.xlsx - add two new model run rows:
#   2035_TM161_FBP_Plan_08 (inserted right after the existing 2035_TM161_FBP_Plan_07 row)
#   2050_TM161_FBP_Plan_03 (appended after the existing 2050 rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new 2035 row right after row 183 (2035_TM161_FBP_Plan_07) ---
# Inserting a whole row here copies the formatting from the row above it,
# which mirrors the source row's (row 183) styling onto the new row 184.
$ws.Rows.Item(184).Insert()

$ws.Range("A184").Value = 2035
$ws.Range("B184").Value = "2035_TM161_FBP_Plan_08"
$ws.Range("C184").Value = "RTP2025"
$ws.Range("D184").Value = "FBP"
$ws.Range("F184").Value = "Update network (v27) and LU (BAUS_FBP_v03)"
$ws.Range("G184").Value = "BAUS_FBP_v03\2035"
$ws.Range("H184").Value = "PBA50Plus_Final_Blueprint_v39"
$ws.Range("I184").Value = "current"
$ws.Range("J184").Value = "FBP"
$ws.Range("K184").Value = "BlueprintNetworks_v27\net_2035_Blueprint"
$ws.Range("L184").Value = "model3-b"
$ws.Range("M184").Value = "https://app.asana.com/1/11860278793487/project/1204085012544660/task/1209660380187449"
$ws.Range("N184").Value = 15.66
$ws.Range("O184").Value = "na"
$ws.Range("P184").Value = "na"
$ws.Range("T184").Value = -0.455
$ws.Range("U184").Value = 5
$ws.Range("V184").Value = 34
$ws.Range("W184").Value = 0
$ws.Range("X184").Value = 85
$ws.Range("Y184").Value = "2035 Plan"

# --- Append the new 2050 row after the last existing row (now row 212) ---
# Inserting here copies formatting from the row above it (row 212, formerly
# row 211, 2050_TM161_FBP_Plan_02) onto the new row 213.
$ws.Rows.Item(213).Insert()

$ws.Range("A213").Value = 2050
$ws.Range("B213").Value = "2050_TM161_FBP_Plan_03"
$ws.Range("C213").Value = "RTP2025"
$ws.Range("D213").Value = "FBP"
$ws.Range("F213").Value = "Update network (v27) and LU (BAUS_FBP_v03)"
$ws.Range("G213").Value = "BAUS_FBP_v03\2050"
$ws.Range("H213").Value = "PBA50Plus_Final_Blueprint_v39"
$ws.Range("I213").Value = "current"
$ws.Range("J213").Value = "FBP"
$ws.Range("K213").Value = "BlueprintNetworks_v27\net_2050_Blueprint"
$ws.Range("L213").Value = "model3-c"
$ws.Range("M213").Value = "https://app.asana.com/1/11860278793487/project/1204085012544660/task/1209660380187451"
$ws.Range("N213").Value = 19.57
$ws.Range("O213").Value = "na"
$ws.Range("P213").Value = "na"
$ws.Range("T213").Value = -0.455
$ws.Range("U213").Value = 5
$ws.Range("V213").Value = 34
$ws.Range("W213").Value = 0
$ws.Range("X213").Value = 85
$ws.Range("Y213").Value = "2050 Plan"
